$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# New rows (TEXT ID, TYPOGRAPHY NAME, ALIGNMENT, DIRECTION, value) for the
# "Grequency mode" panel input / clock UI texts.
$data = @(
    @("SingleUseId55", "Default", "Left",   "LTR", "THRESHOLD"),
    @("SingleUseId56", "Default", "Left",   "LTR", "SLOPE"),
    @("SingleUseId57", "Default", "Center", "LTR", "Detect"),
    @("SingleUseId58", "Default", "Center", "LTR", "<value> mV"),
    @("SingleUseId59", "Default", "Left",   "LTR", "0"),
    @("SingleUseId60", "Default", "Left",   "LTR", "Manual"),
    @("SingleUseId61", "Default", "Left",   "LTR", "Defined"),
    @("SingleUseId62", "Default", "Left",   "LTR", "External"),
    @("SingleUseId63", "Default", "Center", "LTR", "Internal`nRubid"),
    @("SingleUseId64", "Default", "Center", "LTR", "Internal`nQuartz")
)

$row = 55
foreach ($r in $data) {
    $ws.Cells.Item($row, 2).Value = $r[0]
    $ws.Cells.Item($row, 3).Value = $r[1]
    $ws.Cells.Item($row, 4).Value = $r[2]
    $ws.Cells.Item($row, 5).Value = $r[3]

    $fCell = $ws.Cells.Item($row, 6)
    if ($r[4] -eq "0") {
        # Force numeric-looking text to stay text (matches existing "0" cells
        # used elsewhere in this column), then drop the style override so the
        # cell keeps the default (unstyled) formatting.
        $fCell.NumberFormat = "@"
        $fCell.Value = $r[4]
        $fCell.Style = "Normal"
    } else {
        $fCell.Value = $r[4]
    }

    # Avoid leaving a stray explicit row height behind for rows whose text
    # contains an embedded newline (would otherwise mark customHeight="1").
    $ws.Rows.Item($row).AutoFit() | Out-Null

    $row++
}
